$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 2.111846333333333
$ws.Range("H2").Value2 = 6.335539
$ws.Range("I2").Value2 = 0.01909882549924913
$ws.Range("J2").Value2 = 0.01909882549924913
$ws.Range("M2").Value2 = 15.03663066666667
$ws.Range("N2").Value2 = 45.109892
$ws.Range("O2").Value2 = 0.279146411176606
$ws.Range("P2").Value2 = 0.279146411176606
$ws.Range("Q2").Value2 = 31.75505333908756
$ws.Range("R2").Value2 = 285.795480051788
$ws.Range("S2").Value2 = 0.005331368595803647
$ws.Range("T2").Value2 = 0.005331368595803646
$ws.Range("G3").Value2 = 2.111846333333333
$ws.Range("H3").Value2 = 6.335539
$ws.Range("I3").Value2 = 0.01909882549924913
$ws.Range("J3").Value2 = 0.01909882549924913
$ws.Range("O3").Value2 = 0.2673306493381863
$ws.Range("P3").Value2 = 0.2673306493381863
$ws.Range("Q3").Value2 = 30.41091946382311
$ws.Range("R3").Value2 = 273.698275174408
$ws.Range("S3").Value2 = 0.005105701422310981
$ws.Range("T3").Value2 = 0.005105701422310981
$ws.Range("G4").Value2 = 2.111846333333333
$ws.Range("H4").Value2 = 6.335539
$ws.Range("I4").Value2 = 0.01909882549924913
$ws.Range("J4").Value2 = 0.01909882549924913
$ws.Range("M4").Value2 = 22.16851266666667
$ws.Range("N4").Value2 = 66.505538
$ws.Range("O4").Value2 = 0.411545703901694
$ws.Range("P4").Value2 = 0.411545703901694
$ws.Range("Q4").Value2 = 46.81649219055356
$ws.Range("R4").Value2 = 421.348429714982
$ws.Range("S4").Value2 = 0.007860039583784109
$ws.Range("T4").Value2 = 0.007860039583784107
$ws.Range("G5").Value2 = 2.111846333333333
$ws.Range("H5").Value2 = 6.335539
$ws.Range("I5").Value2 = 0.01909882549924913
$ws.Range("J5").Value2 = 0.01909882549924913
$ws.Range("M5").Value2 = 2.261165333333333
$ws.Range("N5").Value2 = 6.783496
$ws.Range("O5").Value2 = 0.04197723558351375
$ws.Range("P5").Value2 = 0.04197723558351374
$ws.Range("Q5").Value2 = 4.775233718260444
$ws.Range("R5").Value2 = 42.977103464344
$ws.Range("S5").Value2 = 0.0008017158973504005
$ws.Range("T5").Value2 = 0.0008017158973504003
$ws.Range("G6").Value2 = 51.92481233333333
$ws.Range("I6").Value2 = 0.4695904783329055
$ws.Range("J6").Value2 = 0.4695904783329055
$ws.Range("M6").Value2 = 15.03663066666667
$ws.Range("N6").Value2 = 45.109892
$ws.Range("O6").Value2 = 0.279146411176606
$ws.Range("P6").Value2 = 0.279146411176606
$ws.Range("Q6").Value2 = 780.7742254923115
$ws.Range("R6").Value2 = 7026.968029430804
$ws.Range("S6").Value2 = 0.1310844967493363
$ws.Range("T6").Value2 = 0.1310844967493363
$ws.Range("G7").Value2 = 51.92481233333333
$ws.Range("I7").Value2 = 0.4695904783329055
$ws.Range("J7").Value2 = 0.4695904783329055
$ws.Range("O7").Value2 = 0.2673306493381863
$ws.Range("P7").Value2 = 0.2673306493381863
$ws.Range("Q7").Value2 = 747.7254671038071
$ws.Range("R7").Value2 = 6729.529203934264
$ws.Range("S7").Value2 = 0.1255359274957651
$ws.Range("T7").Value2 = 0.1255359274957651
$ws.Range("G8").Value2 = 51.92481233333333
$ws.Range("I8").Value2 = 0.4695904783329055
$ws.Range("J8").Value2 = 0.4695904783329055
$ws.Range("M8").Value2 = 22.16851266666667
$ws.Range("N8").Value2 = 66.505538
$ws.Range("O8").Value2 = 0.411545703901694
$ws.Range("P8").Value2 = 0.411545703901694
$ws.Range("Q8").Value2 = 1151.095859925789
$ws.Range("R8").Value2 = 10359.8627393321
$ws.Range("S8").Value2 = 0.1932579439510488
$ws.Range("T8").Value2 = 0.1932579439510488
$ws.Range("G9").Value2 = 51.92481233333333
$ws.Range("I9").Value2 = 0.4695904783329055
$ws.Range("J9").Value2 = 0.4695904783329055
$ws.Range("M9").Value2 = 2.261165333333333
$ws.Range("N9").Value2 = 6.783496
$ws.Range("O9").Value2 = 0.04197723558351375
$ws.Range("P9").Value2 = 0.04197723558351374
$ws.Range("Q9").Value2 = 117.4105855879724
$ws.Range("R9").Value2 = 1056.695270291752
$ws.Range("S9").Value2 = 0.01971211013675528
$ws.Range("T9").Value2 = 0.01971211013675528
$ws.Range("G10").Value2 = 56.36634066666667
$ws.Range("H10").Value2 = 169.099022
$ws.Range("I10").Value2 = 0.5097581615820991
$ws.Range("J10").Value2 = 0.5097581615820991
$ws.Range("M10").Value2 = 15.03663066666667
$ws.Range("N10").Value2 = 45.109892
$ws.Range("O10").Value2 = 0.279146411176606
$ws.Range("P10").Value2 = 0.279146411176606
$ws.Range("Q10").Value2 = 847.5598466361805
$ws.Range("R10").Value2 = 7628.038619725624
$ws.Range("S10").Value2 = 0.1422971613736274
$ws.Range("T10").Value2 = 0.1422971613736274
$ws.Range("G11").Value2 = 56.36634066666667
$ws.Range("H11").Value2 = 169.099022
$ws.Range("I11").Value2 = 0.5097581615820991
$ws.Range("J11").Value2 = 0.5097581615820991
$ws.Range("O11").Value2 = 0.2673306493381863
$ws.Range("P11").Value2 = 0.2673306493381863
$ws.Range("Q11").Value2 = 811.684173904265
$ws.Range("R11").Value2 = 7305.157565138385
$ws.Range("S11").Value2 = 0.1362739803411826
$ws.Range("T11").Value2 = 0.1362739803411826
$ws.Range("G12").Value2 = 56.36634066666667
$ws.Range("H12").Value2 = 169.099022
$ws.Range("I12").Value2 = 0.5097581615820991
$ws.Range("J12").Value2 = 0.5097581615820991
$ws.Range("M12").Value2 = 22.16851266666667
$ws.Range("N12").Value2 = 66.505538
$ws.Range("O12").Value2 = 0.411545703901694
$ws.Range("P12").Value2 = 0.411545703901694
$ws.Range("Q12").Value2 = 1249.557937042649
$ws.Range("R12").Value2 = 11246.02143338384
$ws.Range("S12").Value2 = 0.2097887814279384
$ws.Range("T12").Value2 = 0.2097887814279384
$ws.Range("G13").Value2 = 56.36634066666667
$ws.Range("H13").Value2 = 169.099022
$ws.Range("I13").Value2 = 0.5097581615820991
$ws.Range("J13").Value2 = 0.5097581615820991
$ws.Range("M13").Value2 = 2.261165333333333
$ws.Range("N13").Value2 = 6.783496
$ws.Range("O13").Value2 = 0.04197723558351375
$ws.Range("P13").Value2 = 0.04197723558351374
$ws.Range("Q13").Value2 = 127.4536154823235
$ws.Range("R13").Value2 = 1147.082539340912
$ws.Range("S13").Value2 = 0.02139823843935064
$ws.Range("T13").Value2 = 0.02139823843935064
$ws.Range("E14").Value2 = 3
$ws.Range("F14").Value2 = 1
$ws.Range("G14").Value2 = 0.171671
$ws.Range("H14").Value2 = 0.5150129999999999
$ws.Range("I14").Value2 = 0.001552534585746342
$ws.Range("J14").Value2 = 0.001552534585746342
$ws.Range("M14").Value2 = 15.03663066666667
$ws.Range("N14").Value2 = 45.109892
$ws.Range("O14").Value2 = 0.279146411176606
$ws.Range("P14").Value2 = 0.279146411176606
$ws.Range("Q14").Value2 = 2.581353423177333
$ws.Range("R14").Value2 = 23.232180808596
$ws.Range("S14").Value2 = 0.0004333844578386501
$ws.Range("T14").Value2 = 0.0004333844578386501
$ws.Range("E15").Value2 = 3
$ws.Range("F15").Value2 = 1
$ws.Range("G15").Value2 = 0.171671
$ws.Range("H15").Value2 = 0.5150129999999999
$ws.Range("I15").Value2 = 0.001552534585746342
$ws.Range("J15").Value2 = 0.001552534585746342
$ws.Range("O15").Value2 = 0.2673306493381863
$ws.Range("P15").Value2 = 0.2673306493381863
$ws.Range("Q15").Value2 = 2.472089409570667
$ws.Range("R15").Value2 = 22.248804686136
$ws.Range("S15").Value2 = 0.0004150400789275617
$ws.Range("T15").Value2 = 0.0004150400789275617
$ws.Range("E16").Value2 = 3
$ws.Range("F16").Value2 = 1
$ws.Range("G16").Value2 = 0.171671
$ws.Range("H16").Value2 = 0.5150129999999999
$ws.Range("I16").Value2 = 0.001552534585746342
$ws.Range("J16").Value2 = 0.001552534585746342
$ws.Range("M16").Value2 = 22.16851266666667
$ws.Range("N16").Value2 = 66.505538
$ws.Range("O16").Value2 = 0.411545703901694
$ws.Range("P16").Value2 = 0.411545703901694
$ws.Range("Q16").Value2 = 3.805690737999333
$ws.Range("R16").Value2 = 34.251216641994
$ws.Range("S16").Value2 = 0.0006389389389227033
$ws.Range("T16").Value2 = 0.0006389389389227032
$ws.Range("E17").Value2 = 3
$ws.Range("F17").Value2 = 1
$ws.Range("G17").Value2 = 0.171671
$ws.Range("H17").Value2 = 0.5150129999999999
$ws.Range("I17").Value2 = 0.001552534585746342
$ws.Range("J17").Value2 = 0.001552534585746342
$ws.Range("M17").Value2 = 2.261165333333333
$ws.Range("N17").Value2 = 6.783496
$ws.Range("O17").Value2 = 0.04197723558351375
$ws.Range("P17").Value2 = 0.04197723558351374
$ws.Range("Q17").Value2 = 0.3881765139386666
$ws.Range("R17").Value2 = 3.493588625448
$ws.Range("S17").Value2 = 0.00006517111005742712
$ws.Range("T17").Value2 = 0.00006517111005742711
